$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 35539
$ws.Range("J3").Value = 35539
$ws.Range("L3").Value = 35539
$ws.Range("N3").Value = -35767
$ws.Range("H28").Value = 791.1667
$ws.Range("I28").Value = 825.0625
$ws.Range("K28").Value = 825.0625
$ws.Range("M28").Value = -340.0625
$ws.Range("H80").Value = 1284.4667
$ws.Range("J80").Value = 1396
$ws.Range("L80").Value = 4188
$ws.Range("N80").Value = -6184
$ws.Range("H83").Value = 1284.4667
$ws.Range("J83").Value = 1396
$ws.Range("L83").Value = 12564
$ws.Range("N83").Value = -22548
$ws.Range("H92").Value = 852.25
$ws.Range("J92").Value = 1000
$ws.Range("L92").Value = 1000
$ws.Range("N92").Value = -3496
$ws.Range("H102").Value = 35539
$ws.Range("J102").Value = 35539
$ws.Range("L102").Value = 35539
$ws.Range("N102").Value = -42029
$ws.Range("H129").Value = 1263.3572
$ws.Range("I129").Value = 970.4
$ws.Range("K129").Value = 2911.2
$ws.Range("M129").Value = 2088.8
$ws.Range("H138").Value = 1760.7742
$ws.Range("J138").Value = 2480.5454
$ws.Range("L138").Value = 7441.6362
$ws.Range("N138").Value = -17721.6362

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H101").Value = 40602
$ws.Range("J101").Value = 40602
$ws.Range("L101").Value = 40602
$ws.Range("N101").Value = -47092
$ws.Range("H133").Value = 100000
$ws.Range("J133").Value = 100000
$ws.Range("L133").Value = 100000
$ws.Range("N133").Value = -105060
$ws.Range("H135").Value = 60666.332
$ws.Range("J135").Value = 60666.332
$ws.Range("L135").Value = 60666.332
$ws.Range("N135").Value = -70806.33199999999
$ws.Range("H139").Value = 117665.75
$ws.Range("J139").Value = 117665.75
$ws.Range("L139").Value = 117665.75
$ws.Range("N139").Value = -127945.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3451.1667
$ws.Range("I86").Value = 3081
$ws.Range("J86").Value = 3715.5715
$ws.Range("K86").Value = 3081
$ws.Range("L86").Value = 3715.5715
$ws.Range("M86").Value = -1958
$ws.Range("N86").Value = -5961.5715
$ws.Range("H89").Value = 3451.1667
$ws.Range("I89").Value = 3081
$ws.Range("J89").Value = 3715.5715
$ws.Range("K89").Value = 15405
$ws.Range("L89").Value = 18577.8575
$ws.Range("M89").Value = -9789
$ws.Range("N89").Value = -29809.8575

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1254.8572
$ws.Range("I31").Value = 1254.8572
$ws.Range("K31").Value = 1254.8572
$ws.Range("M31").Value = -959.8571999999999
$ws.Range("H34").Value = 1254.8572
$ws.Range("I34").Value = 1254.8572
$ws.Range("K34").Value = 1254.8572
$ws.Range("M34").Value = -1052.8572
$ws.Range("H51").Value = 9999.714
$ws.Range("J51").Value = 9999.666999999999
$ws.Range("L51").Value = 9999.666999999999
$ws.Range("N51").Value = -11471.667
$ws.Range("H60").Value = 10006.143
$ws.Range("J60").Value = 9999.5
$ws.Range("L60").Value = 9999.5
$ws.Range("N60").Value = -11021.5
$ws.Range("H61").Value = 9999.714
$ws.Range("J61").Value = 9999.666999999999
$ws.Range("L61").Value = 9999.666999999999
$ws.Range("N61").Value = -10695.667
$ws.Range("H107").Value = 1570.4117
$ws.Range("I107").Value = 920.7143
$ws.Range("K107").Value = 920.7143
$ws.Range("M107").Value = 999.2857
$ws.Range("H122").Value = 3204.1765
$ws.Range("I122").Value = 4358
$ws.Range("J122").Value = 1906.125
$ws.Range("K122").Value = 13074
$ws.Range("L122").Value = 5718.375
$ws.Range("M122").Value = -10624
$ws.Range("N122").Value = -10618.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2005
$ws.Range("I3").Value = 2005
$ws.Range("K3").Value = 6015
$ws.Range("M3").Value = -5903
$ws.Range("H62").Value = 206799
$ws.Range("J62").Value = 8999.5
$ws.Range("L62").Value = 26998.5
$ws.Range("N62").Value = -28370.5
$ws.Range("H63").Value = 11999.5
$ws.Range("J63").Value = 19999
$ws.Range("L63").Value = 59997
$ws.Range("N63").Value = -61495
$ws.Range("H65").Value = 206799
$ws.Range("J65").Value = 8999.5
$ws.Range("L65").Value = 80995.5
$ws.Range("N65").Value = -87859.5
$ws.Range("H66").Value = 11999.5
$ws.Range("J66").Value = 19999
$ws.Range("L66").Value = 179991
$ws.Range("N66").Value = -187479
$ws.Range("H81").Value = 71430840
$ws.Range("I81").Value = 2597.4
$ws.Range("J81").Value = 250001460
$ws.Range("K81").Value = 7792.200000000001
$ws.Range("L81").Value = 750004380
$ws.Range("M81").Value = -6669.200000000001
$ws.Range("N81").Value = -750006626
$ws.Range("H84").Value = 71430840
$ws.Range("I84").Value = 2597.4
$ws.Range("J84").Value = 250001460
$ws.Range("K84").Value = 23376.6
$ws.Range("L84").Value = 2250013140
$ws.Range("M84").Value = -17760.6
$ws.Range("N84").Value = -2250024372
$ws.Range("H107").Value = 666
$ws.Range("I107").Value = 363.53845
$ws.Range("J107").Value = 801.5862
$ws.Range("K107").Value = 1090.61535
$ws.Range("L107").Value = 2404.7586
$ws.Range("M107").Value = 829.38465
$ws.Range("N107").Value = -6244.7586
$ws.Range("H114").Value = 20001248
$ws.Range("J114").Value = 2249.25
$ws.Range("L114").Value = 6747.75
$ws.Range("N114").Value = -13255.75
$ws.Range("H129").Value = 113543.555
$ws.Range("J129").Value = 4553.4443
$ws.Range("L129").Value = 13660.3329
$ws.Range("N129").Value = -23660.3329
$ws.Range("H137").Value = 6671265
$ws.Range("I137").Value = 12502205
$ws.Range("K137").Value = 37506615
$ws.Range("M137").Value = -37501515

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8807.23
$ws.Range("I7").Value = 8283.666999999999
$ws.Range("J7").Value = 9256
$ws.Range("K7").Value = 8283.666999999999
$ws.Range("L7").Value = 9256
$ws.Range("M7").Value = -8171.666999999999
$ws.Range("N7").Value = -9480
$ws.Range("H82").Value = 1481.6364
$ws.Range("J82").Value = 1618.1666
$ws.Range("L82").Value = 1618.1666
$ws.Range("N82").Value = -2340.1666
$ws.Range("H85").Value = 1481.6364
$ws.Range("J85").Value = 1618.1666
$ws.Range("L85").Value = 1618.1666
$ws.Range("N85").Value = -4114.1666
$ws.Range("H93").Value = 1493.5
$ws.Range("I93").Value = 1255.6316
$ws.Range("K93").Value = 1255.6316
$ws.Range("M93").Value = -7.631599999999935
$ws.Range("H100").Value = 9433
$ws.Range("I100").Value = 18499
$ws.Range("J100").Value = 4900
$ws.Range("K100").Value = 18499
$ws.Range("L100").Value = 4900
$ws.Range("M100").Value = -17958
$ws.Range("N100").Value = -5982
$ws.Range("H121").Value = 75485
$ws.Range("J121").Value = 75485
$ws.Range("L121").Value = 75485
$ws.Range("N121").Value = -78979
$ws.Range("H126").Value = 8807.23
$ws.Range("I126").Value = 8283.666999999999
$ws.Range("J126").Value = 9256
$ws.Range("K126").Value = 24851.001
$ws.Range("L126").Value = 27768
$ws.Range("M126").Value = -22381.001
$ws.Range("N126").Value = -32708

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6214.8237
$ws.Range("J62").Value = 4830.5
$ws.Range("L62").Value = 4830.5
$ws.Range("N62").Value = -6078.5
$ws.Range("H65").Value = 6214.8237
$ws.Range("J65").Value = 4830.5
$ws.Range("L65").Value = 24152.5
$ws.Range("N65").Value = -30392.5
$ws.Range("H105").Value = 30307.5
$ws.Range("J105").Value = 30307.5
$ws.Range("L105").Value = 30307.5
$ws.Range("N105").Value = -37295.5
$ws.Range("H122").Value = 2257.5386
$ws.Range("I122").Value = 300.66666
$ws.Range("K122").Value = 901.9999799999999
$ws.Range("M122").Value = 1548.00002
$ws.Range("H136").Value = 761.0625
$ws.Range("I136").Value = 740.13794
$ws.Range("J136").Value = 963.3333
$ws.Range("K136").Value = 2220.41382
$ws.Range("L136").Value = 2889.9999
$ws.Range("M136").Value = 329.5861800000002
$ws.Range("N136").Value = -7989.9999
